# Fix notification import user and add transactional insert database
#
# The "userImport.xlsx" sample sheet gains two new sample rows (new rows 2
# and 3) inserted above the previous two data rows (which shift down to
# become rows 4 and 5). The D column for the new second row is entered as
# text ("asdf") instead of a number, matching the authored fixture data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the existing data rows 2:3, pushing the
# existing sample rows down to rows 4:5.
$ws.Range("A2:G3").Insert()

# New row 2 sample data
$ws.Range("A2").Value = "aldo"
$ws.Range("B2").Value = "asd"
$ws.Range("C2").Value = "asdasd"
$ws.Range("D2").Value = "asdf"
$ws.Range("E2").Value = "anasdkajsdh"
$ws.Range("F2").Value = "asdf@sad"
$ws.Range("G2").Value = "admin"

# New row 3 sample data
$ws.Range("A3").Value = "andrianasd2"
$ws.Range("B3").Value = "andriasdn2"
$ws.Range("C3").Value = "dasdf"
$ws.Range("D3").Value = 9871
$ws.Range("E3").Value = "dasd"
$ws.Range("F3").Value = "adaff@ad"
$ws.Range("G3").Value = "lecturer"

# Match the selection recorded in the saved workbook (active cell F3).
$ws.Range("F3").Select()
